# Peer-grading "forms_responses" sheet: correct the date on the existing
# test-case row and append the real survey rows received from the form,
# then grow the table/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Named cell style used for the RaterEmail/RaterName columns whenever
#     the rater is "ajaysathish.shenoy@uzh.ch" or "benjaminlucasde.gorgey@uzh.ch"
#     (font size bumped from 11 to 12, same Calibri/theme color otherwise). ---
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.Size = 12

# --- Row 3 (existing test row): fix the submitter e-mail/name pair and the
#     ReceivedAtUTC timestamp (typo'd month corrected from March to February). ---
$ws.Cells.Item(3, 3).Value = "ajaysathish.shenoy@uzh.ch"
$ws.Cells.Item(3, 3).Style = "Normal 2"
$ws.Cells.Item(3, 4).Value = "ajaysathish.shenoy@uzh.ch"
$ws.Cells.Item(3, 4).Style = "Normal 2"
$ws.Cells.Item(3, 14).Value = "2026-02-19T14:35:52.5219332Z"
$ws.Rows.Item(3).RowHeight = 15.6

# --- New survey responses appended as rows 4-8. ---
$nbsp = [char]0x00A0
$data = @(
    @{ Row = 4;  ResponseId = 3; Rater = "arda.aydin@uzh.ch";             Styled = $false;
       Name = "Abirami`tSivarajah"; Q2_1 = "Good 4"; Q2_2 = "Excellent 5"; Q2_3 = "Fair 3";
       Q3_1 = "Fair 3"; Q3_3 = "Good 4"; Q3_2 = "Poor 2"; Q_4 = "Poor 2";
       Comment = "testcase"; ReceivedAtUTC = "2026-02-20T09:50:20.4657318Z" }

    @{ Row = 5;  ResponseId = 4; Rater = "arda.aydin@uzh.ch";             Styled = $false;
       Name = "Gizem`tTopsakal"; Q2_1 = "$nbsp" + "Bad 1"; Q2_2 = "Excellent 5"; Q2_3 = "Excellent 5";
       Q3_1 = "Excellent 5"; Q3_3 = "Poor 2"; Q3_2 = "Bad 1"; Q_4 = "Excellent 5";
       Comment = "testcase"; ReceivedAtUTC = "2026-02-20T09:50:50.0443991Z" }

    @{ Row = 6;  ResponseId = 5; Rater = "ajaysathish.shenoy@uzh.ch";     Styled = $true;
       Name = "Anna-Lea`tWölfle"; Q2_1 = "Fair 3"; Q2_2 = "Good 4"; Q2_3 = "Poor 2";
       Q3_1 = "Poor 2"; Q3_3 = "Fair 3"; Q3_2 = "Poor 2"; Q_4 = "Excellent 5";
       Comment = "testcase"; ReceivedAtUTC = "2026-02-20T09:51:05.1273717Z" }

    @{ Row = 7;  ResponseId = 6; Rater = "benjaminlucasde.gorgey@uzh.ch"; Styled = $true;
       Name = "Giacomo`tBasile"; Q2_1 = "Excellent 5"; Q2_2 = "Excellent 5"; Q2_3 = "Excellent 5";
       Q3_1 = "Excellent 5"; Q3_3 = "Excellent 5"; Q3_2 = "Excellent 5"; Q_4 = "Excellent 5";
       Comment = "testcase"; ReceivedAtUTC = "2026-02-20T09:51:20.7604177Z" }

    @{ Row = 8;  ResponseId = 7; Rater = "ajaysathish.shenoy@uzh.ch";     Styled = $true;
       Name = "Alejandro`tCaicedo Murgueitio"; Q2_1 = "Poor 2"; Q2_2 = "Fair 3"; Q2_3 = "Good 4";
       Q3_1 = "Good 4"; Q3_3 = "Poor 2"; Q3_2 = "Fair 3"; Q_4 = "Good 4";
       Comment = "tescase"; ReceivedAtUTC = "2026-02-20T09:58:26.8961547Z" }
)

foreach ($rec in $data) {
    $r = $rec.Row

    $ws.Cells.Item($r, 1).Value = $rec.ResponseId   # A: ResponseId
    # B (SubmittedAt) intentionally left blank, matching the source rows.

    $ws.Cells.Item($r, 3).Value = $rec.Rater
    $ws.Cells.Item($r, 4).Value = $rec.Rater
    if ($rec.Styled) {
        $ws.Cells.Item($r, 3).Style = "Normal 2"
        $ws.Cells.Item($r, 4).Style = "Normal 2"
        $ws.Rows.Item($r).RowHeight = 15.6
    }

    $ws.Cells.Item($r, 5).Value  = $rec.Name
    $ws.Cells.Item($r, 6).Value  = $rec.Q2_1
    $ws.Cells.Item($r, 7).Value  = $rec.Q2_2
    $ws.Cells.Item($r, 8).Value  = $rec.Q2_3
    $ws.Cells.Item($r, 9).Value  = $rec.Q3_1
    $ws.Cells.Item($r, 10).Value = $rec.Q3_3
    $ws.Cells.Item($r, 11).Value = $rec.Q3_2
    $ws.Cells.Item($r, 12).Value = $rec.Q_4
    $ws.Cells.Item($r, 13).Value = $rec.Comment
    $ws.Cells.Item($r, 14).Value = $rec.ReceivedAtUTC
}

# --- Grow the "Responses" table + AutoFilter to cover the new rows. ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:N8"))

# --- Sheet dimension / selection as left by the editor. ---
$ws.Range("C7").Select()

Write-Output "done"
